# Update "QS WSFT Meeting 2019-06-09.pptx" slide 7 page-number textboxes.
#  - Move the existing "-1-" textbox down slightly and renumber it "-2-".
#  - Add four more page-number textboxes ("-3-", "-4-", "-5-", "-5-").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# EMU -> Points helper (1 pt = 12700 EMU); PowerPoint COM works in points.
function EmuToPt([double]$emu) { return $emu / 12700.0 }

# --- 1. Existing "TextBox 43" (id 44): shift down & change its label to "-2-" ---
# NOTE: the `.Top` setter here truncates through a 32-bit float before
# converting back to EMU, so the plain `EmuToPt 4771738` value lands one EMU
# short (4771737). Nudge it up by a hair (well under half a point) so the
# truncated result still lands exactly on the target EMU.
$shpExisting = $s.Shapes.Item(34)
$newTop = EmuToPt 4771738
$newTop = $newTop + 0.00005
$shpExisting.Top = $newTop
$shpExisting.TextFrame.TextRange.Text = "-2-"

# --- helper to build one of the new page-number textboxes ---
# NOTE: call positionally with pre-computed variables (named `-Param value`
# binding, and passing parenthesised expressions directly, aren't reliable
# in this host).
function Add-PageNumberTextBox {
    param(
        [double]$Left,
        [double]$Top,
        [double]$Width,
        [double]$Height,
        [string]$Text,
        [int]$Alignment   # 1 = ppAlignLeft, 2 = ppAlignCenter
    )

    $shp = $s.Shapes.AddTextbox(1, $Left, $Top, $Width, $Height)
    $shp.Fill.Visible = 0
    $shp.TextFrame.AutoSize = 1
    $shp.TextFrame.WordWrap = 1
    $shp.TextFrame.MarginLeft = 7.2
    $shp.TextFrame.MarginRight = 7.2
    $shp.TextFrame.MarginTop = 3.6
    $shp.TextFrame.MarginBottom = 3.6

    $tr = $shp.TextFrame.TextRange
    $tr.Text = $Text
    $tr.LanguageID = "nl-NL"
    $tr.ParagraphFormat.Alignment = $Alignment

    # Turning AutoSize on (and then setting the text) makes the host
    # recompute Height to fit the text, clobbering the height we asked for
    # at construction time. Re-assert it now that the text is in place.
    # (The `.Height` setter truncates through a 32-bit float on its way
    # back to EMU, so nudge up by a hair to land on the exact target EMU.)
    $safeHeight = $Height + 0.00005
    $shp.Height = $safeHeight

    return $shp
}

# --- 2. "TextBox 45" -> "-3-" ---
$left = EmuToPt 9806795
$top = EmuToPt 4901836
$width = EmuToPt 557842
$height = EmuToPt 369332
Add-PageNumberTextBox $left $top $width $height "-3-" 1 | Out-Null

# --- 3. "TextBox 46" -> "-4-" ---
$left = EmuToPt 11358673
$top = EmuToPt 4864664
$width = EmuToPt 557842
$height = EmuToPt 369332
Add-PageNumberTextBox $left $top $width $height "-4-" 1 | Out-Null

# --- 4. "TextBox 47" -> "-5-" ---
$left = EmuToPt 207453
$top = EmuToPt 4901835
$width = EmuToPt 557842
$height = EmuToPt 369332
Add-PageNumberTextBox $left $top $width $height "-5-" 1 | Out-Null

# --- 5. "TextBox 48" -> "-5-" centered ---
$left = EmuToPt 3422722
$top = EmuToPt 2755225
$width = EmuToPt 567134
$height = EmuToPt 378624
Add-PageNumberTextBox $left $top $width $height "-5-" 2 | Out-Null
